$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the existing "sum" header cell (G1) onto the new H1 header
# so that H1 gets the same bold/border/alignment style (style index 1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New "Save" column values for rows 2-17 (plain numeric cells, no special style,
# matching the style of the other numeric data columns).
$saveValues = @(0, 0, 1, 1, 0, 1, 0, 0, 0, 1, 1, 1, 0, 1, 1, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
